$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D; existing D:K data shifts to E:L
$ws.Columns("D:D").Insert()

# New column D should carry the same number formatting as the column
# immediately to its right (the old column D, now shifted to E)
$ws.Range("E5:E102").Copy() | Out-Null
$ws.Range("D5:D102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate the new column D with the latest reporting-period figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 2973900
$ws.Range("D9").Value = 433800
$ws.Range("D10").Value = 2540100
$ws.Range("D12").Value = 440000
$ws.Range("D14").Value = 21300
$ws.Range("D15").Value = 15900
$ws.Range("D17").Value = 2300500
$ws.Range("D18").Value = 673400
$ws.Range("D20").Value = 36300
$ws.Range("D21").Value = 851500
$ws.Range("D22").Value = 80200
$ws.Range("D23").Value = 629500
$ws.Range("D24").Value = 80100
$ws.Range("D26").Value = 549400
$ws.Range("D27").Value = 549400
$ws.Range("D29").Value = 26300
$ws.Range("D32").Value = -36300
$ws.Range("D33").Value = 575700
$ws.Range("D35").Value = 575700
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 618800
$ws.Range("D42").Value = 583600
$ws.Range("D43").Value = 688400
$ws.Range("D44").Value = 21900
$ws.Range("D45").Value = 174200
$ws.Range("D46").Value = 2086900
$ws.Range("D47").Value = 574300
$ws.Range("D48").Value = 243400
$ws.Range("D49").Value = 1969900
$ws.Range("D52").Value = 261600
$ws.Range("D54").Value = 5136000
$ws.Range("D57").Value = 75600
$ws.Range("D59").Value = 1680100
$ws.Range("D60").Value = 2911100
$ws.Range("D61").Value = 749900
$ws.Range("D62").Value = 923500
$ws.Range("D66").Value = 4584500
$ws.Range("D72").Value = 4169000
$ws.Range("D76").Value = 551500
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 575700
$ws.Range("D83").Value = 141900
$ws.Range("D89").Value = 1035300
$ws.Range("D94").Value = 132200
$ws.Range("D100").Value = -1658000
$ws.Range("D101").Value = -5800
$ws.Range("D102").Value = -496400

# Row 91 (Changes In Other Operating Activities) was recalculated for every
# period -- not a simple shift -- so overwrite D91:K91 explicitly
$ws.Range("D91").Value = -69400
$ws.Range("E91").Value = -80900
$ws.Range("F91").Value = -85000
$ws.Range("G91").Value = -97600
$ws.Range("H91").Value = -165400
$ws.Range("I91").Value = -162900
$ws.Range("J91").Value = -123000
$ws.Range("K91").Value = -144200
